# Applies the two changes from the commit:
#  1. Slide 5's table switches from the deck's custom table style
#     ({E0993017-878E-4ABC-B0FF-246345ED91A3}) to the built-in table
#     style {7387B073-A070-458D-8CD6-9ED75AF9A958}.
#  2. The presentation's applied theme colour scheme changes from the
#     "Integral" ("Red Violet") palette to the default "Office Theme"
#     palette (dk1/lt1/dk2/lt2/accent1-6/hlink/folHlink).

$p = $ppt.ActivePresentation

# --- 1. Table style -------------------------------------------------
$tableSlide = $p.Slides.Item(5)
$tableShape = $tableSlide.Shapes.Item(2)
$tbl = $tableShape.Table
$tbl.ApplyStyle("{7387B073-A070-458D-8CD6-9ED75AF9A958}")

# --- 2. Theme colour scheme ------------------------------------------
# The Theme colour scheme is shared across the whole deck, so it can be
# reached from any slide; use slide 1. RGB values below are encoded as
# PowerPoint's COM "long" (0xBBGGRR, i.e. RGB(r,g,b) = r + g*256 + b*65536)
# matching the standard Office theme palette.
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme
$tcs.Item(1).RGB  = 0x000000   # dk1      = 000000
$tcs.Item(2).RGB  = 0xFFFFFF   # lt1      = FFFFFF
$tcs.Item(3).RGB  = 0x6A5444   # dk2      = 44546A
$tcs.Item(4).RGB  = 0xE6E6E7   # lt2      = E7E6E6
$tcs.Item(5).RGB  = 0xD59B5B   # accent1  = 5B9BD5
$tcs.Item(6).RGB  = 0x317DED   # accent2  = ED7D31
$tcs.Item(7).RGB  = 0xA5A5A5   # accent3  = A5A5A5
$tcs.Item(8).RGB  = 0x00C0FF   # accent4  = FFC000
$tcs.Item(9).RGB  = 0xC47244   # accent5  = 4472C4
$tcs.Item(10).RGB = 0x47AD70   # accent6  = 70AD47
$tcs.Item(11).RGB = 0xC16305   # hlink    = 0563C1
$tcs.Item(12).RGB = 0x724F95   # folHlink = 954F72
